# Apply cryptos list price/volume updates (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.906.14"
$ws.Range("E2").Value = "  +3.62%  "
$ws.Range("D3").Value = "3.052.31"
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'557.15"
$ws.Range("E5").Value = "  +3.02%  "
$ws.Range("D6").Value = "'142.17"
$ws.Range("E6").Value = "  +6.18%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.050.58"
$ws.Range("E8").Value = "  +2.70%  "
$ws.Range("D9").Value = "'0.507"
$ws.Range("E9").Value = "  +4.67%  "
$ws.Range("E10").Value = "  +6.92%  "
$ws.Range("D11").Value = "'6.06"
$ws.Range("E11").Value = "  -9.27%  "
$ws.Range("D12").Value = "'0.476"
$ws.Range("E12").Value = "  +7.95%  "
$ws.Range("D13").Value = "'0.0000230"
$ws.Range("E13").Value = "  +6.46%  "
$ws.Range("D14").Value = "'34.89"
$ws.Range("E14").Value = "  +4.68%  "
$ws.Range("D15").Value = "3.550.85"
$ws.Range("E15").Value = "  +4.15%  "
$ws.Range("D16").Value = "63.932.40"
$ws.Range("E16").Value = "  +3.63%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.054.99"
$ws.Range("E17").Value = "  +3.09%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.110"
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("E19").Value = "  +3.27%  "
$ws.Range("D20").Value = "'474.36"
$ws.Range("E20").Value = "  +2.59%  "
$ws.Range("D21").Value = "'14.00"
$ws.Range("E21").Value = "  +5.22%  "
$ws.Range("E22").Value = "  +4.58%  "
$ws.Range("E23").Value = "  +7.29%  "
$ws.Range("D24").Value = "'14.20"
$ws.Range("E24").Value = "  +14.37%  "
$ws.Range("D25").Value = "'81.38"
$ws.Range("E25").Value = "  +3.58%  "
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("D27").Value = "'2.78"
$ws.Range("E27").Value = "  +3.42%  "
$ws.Range("D28").Value = "'7.90"
$ws.Range("E28").Value = "  +5.09%  "
$ws.Range("E29").Value = "  +2.89%  "
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("D31").Value = "'26.20"
$ws.Range("E31").Value = "  +4.85%  "
$ws.Range("E32").Value = "  +1.94%  "
$ws.Range("E33").Value = "  +5.88%  "
$ws.Range("D34").Value = "'5.57"
$ws.Range("E34").Value = "  +2.45%  "
$ws.Range("D35").Value = "'6.18"
$ws.Range("E35").Value = "  +7.55%  "
$ws.Range("D36").Value = "'54.79"
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("D37").Value = "'0.0404"
$ws.Range("E37").Value = "  +5.51%  "
$ws.Range("D38").Value = "'441.88"
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("D39").Value = "'0.0804"
$ws.Range("E39").Value = "  +1.29%  "
$ws.Range("D40").Value = "'2.82"
$ws.Range("E40").Value = "  +17.48%  "
$ws.Range("D41").Value = "2.960.99"
$ws.Range("E41").Value = "  +1.65%  "
$ws.Range("E42").Value = "  +3.59%  "
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("D44").Value = "'27.61"
$ws.Range("E44").Value = "  +4.55%  "
$ws.Range("E45").Value = "  +6.13%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("E47").Value = "  +8.63%  "
$ws.Range("D48").Value = "'0.112"
$ws.Range("E48").Value = "  +4.93%  "
$ws.Range("D49").Value = "'117.08"
$ws.Range("E49").Value = "  +3.32%  "
$ws.Range("D50").Value = "0.0₃0511"
$ws.Range("E50").Value = "  +6.35%  "
$ws.Range("E51").Value = "  +4.31%  "
